$wb = $excel.ActiveWorkbook

# --- Summary sheet updates ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B5").Value = 0.4
$summary.Range("B6").Value = 31
$summary.Range("B9").Value = 45.16

# --- Strategy Status sheet updates (MarketMaking row, row 4) ---
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("D4").Value = 31
$status.Range("G4").Value = 45.16

# --- New trade row (#31) data, shared by "All Trades" and "MarketMaking" sheets ---
$rowValues = @("31", "2026-02-17", "04:15:29", "MarketMaking", "DOWN", 0.01, 0.01, "CLOSED", 0, 0, 100.62, 0, 0, 0.6, "Normal spread capture: 19600 bps", "early_exit", 0.13)

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $newRow = 32
    # Force the date-like column to be stored as literal text, not auto-converted to a date serial
    $ws.Range("B$newRow").NumberFormat = "@"
    for ($col = 1; $col -le $rowValues.Length; $col++) {
        $ws.Cells.Item($newRow, $col).Value = $rowValues[$col - 1]
    }
}
